$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.5544143368907429
$ws.Range("J2").Value = 0.5544143368907429
$ws.Range("M2").Value = 1.443038
$ws.Range("N2").Value = 4.329114
$ws.Range("O2").Value = 0.0289666880885598
$ws.Range("P2").Value = 0.0289666880885598
$ws.Range("Q2").Value = 0.2851159290526666
$ws.Range("R2").Value = 2.566043361474
$ws.Range("S2").Value = 0.01605954716853986
$ws.Range("T2").Value = 0.01605954716853986
$ws.Range("I3").Value = 0.5544143368907429
$ws.Range("J3").Value = 0.5544143368907429
$ws.Range("N3").Value = 87.61054300000001
$ws.Range("O3").Value = 0.5862140087672342
$ws.Range("P3").Value = 0.5862140087672342
$ws.Range("Q3").Value = 5.770040096484778
$ws.Range("R3").Value = 51.930360868363
$ws.Range("S3").Value = 0.3250054509467503
$ws.Range("T3").Value = 0.3250054509467503
$ws.Range("I4").Value = 0.5544143368907429
$ws.Range("J4").Value = 0.5544143368907429
$ws.Range("M4").Value = 19.170603
$ws.Range("N4").Value = 57.511809
$ws.Range("O4").Value = 0.384819303144206
$ws.Range("P4").Value = 0.384819303144206
$ws.Range("Q4").Value = 3.787734130941
$ws.Range("R4").Value = 34.089607178469
$ws.Range("S4").Value = 0.2133493387754528
$ws.Range("T4").Value = 0.2133493387754528
$ws.Range("G5").Value = 0.1587963333333333
$ws.Range("H5").Value = 0.476389
$ws.Range("I5").Value = 0.4455856631092571
$ws.Range("J5").Value = 0.4455856631092571
$ws.Range("M5").Value = 1.443038
$ws.Range("N5").Value = 4.329114
$ws.Range("O5").Value = 0.0289666880885598
$ws.Range("P5").Value = 0.0289666880885598
$ws.Range("Q5").Value = 0.2291491432606667
$ws.Range("R5").Value = 2.062342289346
$ws.Range("S5").Value = 0.01290714092001994
$ws.Range("T5").Value = 0.01290714092001994
$ws.Range("G6").Value = 0.1587963333333333
$ws.Range("H6").Value = 0.476389
$ws.Range("I6").Value = 0.4455856631092571
$ws.Range("J6").Value = 0.4455856631092571
$ws.Range("N6").Value = 87.61054300000001
$ws.Range("O6").Value = 0.5862140087672342
$ws.Range("P6").Value = 0.5862140087672342
$ws.Range("Q6").Value = 4.637410996580778
$ws.Range("R6").Value = 41.73669896922701
$ws.Range("S6").Value = 0.2612085578204839
$ws.Range("T6").Value = 0.2612085578204839
$ws.Range("G7").Value = 0.1587963333333333
$ws.Range("H7").Value = 0.476389
$ws.Range("I7").Value = 0.4455856631092571
$ws.Range("J7").Value = 0.4455856631092571
$ws.Range("M7").Value = 19.170603
$ws.Range("N7").Value = 57.511809
$ws.Range("O7").Value = 0.384819303144206
$ws.Range("P7").Value = 0.384819303144206
$ws.Range("Q7").Value = 3.044221464189
$ws.Range("R7").Value = 27.397993177701
$ws.Range("S7").Value = 0.1714699643687533
$ws.Range("T7").Value = 0.1714699643687533